# Update column F (dSF) values for several rows per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value  = -2
$ws.Range("F4").Value  = 2
$ws.Range("F5").Value  = 1
$ws.Range("F6").Value  = 0
$ws.Range("F8").Value  = 10
$ws.Range("F9").Value  = -5
$ws.Range("F10").Value = 1
$ws.Range("F12").Value = -1
